$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text representation
# (values like "0.9994" or "242.77" would otherwise be auto-converted
# to numbers by Excel, losing formatting such as trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.171.03"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "1.893.45"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "0.7442"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "242.77"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.9994"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "0.3174"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "0.07243"
$ws.Range("D10").Value = "25.04"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("D11").Value = "0.08364"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").Value = "0.7639"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "5.455"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "1.900.89"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "93.03"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "6.186"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "30.179.61"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").Value = "250.88"
$ws.Range("E18").Value = "  +3.10%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "0.000007872"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").Value = "2.159.98"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "0.9990"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "8.018"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "0.9993"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "0.1587"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "9.319"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "164.41"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").Value = "2.077"
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("D30").Value = "1.478"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D31").Value = "4.613"
$ws.Range("D32").Value = "1.537"
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").Value = "4.230"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").Value = "0.05413"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").Value = "1.256"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").Value = "0.7686"
$ws.Range("E36").Value = "  +3.32%  "
$ws.Range("D37").Value = "0.9945"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").Value = "2.720"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  +2.36%  "
$ws.Range("D40").Value = "2.774"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "0.4580"
$ws.Range("E41").Value = "  +2.92%  "
$ws.Range("D42").Value = "1.102.25"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").Value = "6.091"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "72.95"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").Value = "0.8708"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").Value = "104.45"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").Value = "1.001"
$ws.Range("D48").Value = "1.874"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").Value = "7.635"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "9.648"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").Value = "2.053.57"
$ws.Range("E51").Value = "  -0.17%  "
